# fix: revert admin dev default; seed customers only when table empty;
# autosave on customer select when hours/day present
#
# The seeded sample timesheet rows get their auto-generated customer
# (client) names replaced with the real/new customer names, the
# dev-seeded Rate/Total defaults are reverted back to 0, and the
# "Seeded sample hours" note is cleared out (shared string removed),
# while the employee id gets re-generated.

$wb = $excel.ActiveWorkbook

$timesheet = $wb.Worksheets.Item("Weekly Timesheet")
$schema = $wb.Worksheets.Item("Jason Schema")

# --- New client names for the 5 seeded rows (rows 2-6) ---
$clients = @("Cobb", "Evans", "Davis", "Funke", "Field")

for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = 2 + $i

    # Weekly Timesheet sheet: Client (B), Rate (E), Total (F)
    $timesheet.Cells.Item($row, 2).Value = $clients[$i]
    $timesheet.Cells.Item($row, 5).Value = 0
    $timesheet.Cells.Item($row, 6).Value = 0

    # Jason Schema sheet: Client (D), Rate (F), Total (G), Notes (I)
    $schema.Cells.Item($row, 4).Value = $clients[$i]
    $schema.Cells.Item($row, 6).Value = 0
    $schema.Cells.Item($row, 7).Value = 0
    $schema.Cells.Item($row, 9).Value = ""
}

# --- Revert the dev-seeded subtotal / grand total defaults back to 0 ---
$timesheet.Cells.Item(8, 6).Value = 0
$timesheet.Cells.Item(11, 6).Value = 0
$timesheet.Cells.Item(13, 6).Value = 0

# --- Re-generate the employee id (admin dev default reverted) ---
$schema.Range("B2:B6").Value = "emp_0tnwvsb4"
